# Auto-generated edit script applying numeric corrections to the
# per-job "Valefor Profits" leve-crafting tables (ALC/ARM/BSM/CRP/CUL/LTW/WVR).
# Columns: H=currentAveragePrice I=currentAveragePriceNQ J=currentAveragePriceHQ
#          K=LevePriceNQ L=LevePriceHQ M=LeveProfitNQ N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H33").Value = 43488584
$ws.Range("I33").Value = 62502190
$ws.Range("K33").Value = 62502190
$ws.Range("M33").Value = -62501961
$ws.Range("H64").Value = 3277.5334
$ws.Range("I64").Value = 2760
$ws.Range("J64").Value = 3622.5557
$ws.Range("K64").Value = 2760
$ws.Range("L64").Value = 3622.5557
$ws.Range("M64").Value = -2512
$ws.Range("N64").Value = -4118.5557
$ws.Range("H67").Value = 3277.5334
$ws.Range("I67").Value = 2760
$ws.Range("J67").Value = 3622.5557
$ws.Range("K67").Value = 2760
$ws.Range("L67").Value = 3622.5557
$ws.Range("M67").Value = -1902
$ws.Range("N67").Value = -5338.5557
$ws.Range("H80").Value = 602.05884
$ws.Range("I80").Value = 547.75
$ws.Range("J80").Value = 650.3333
$ws.Range("K80").Value = 1643.25
$ws.Range("L80").Value = 1950.9999
$ws.Range("M80").Value = -645.25
$ws.Range("N80").Value = -3946.9999
$ws.Range("H83").Value = 602.05884
$ws.Range("I83").Value = 547.75
$ws.Range("J83").Value = 650.3333
$ws.Range("K83").Value = 4929.75
$ws.Range("L83").Value = 5852.9997
$ws.Range("M83").Value = 62.25
$ws.Range("N83").Value = -15836.9997
$ws.Range("H92").Value = 472.94446
$ws.Range("I92").Value = 202.4
$ws.Range("J92").Value = 811.125
$ws.Range("K92").Value = 202.4
$ws.Range("L92").Value = 811.125
$ws.Range("M92").Value = 1045.6
$ws.Range("N92").Value = -3307.125
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H103").Value = 2778502.8
$ws.Range("J103").Value = 966.6667
$ws.Range("L103").Value = 2900.0001
$ws.Range("N103").Value = -4072.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1093.3334
$ws.Range("I2").Value = 997.7778
$ws.Range("J2").Value = 1666.6666
$ws.Range("K2").Value = 997.7778
$ws.Range("L2").Value = 1666.6666
$ws.Range("M2").Value = -884.7778
$ws.Range("N2").Value = -1892.6666
$ws.Range("H61").Value = 1282.8064
$ws.Range("I61").Value = 1076.4762
$ws.Range("J61").Value = 1716.1
$ws.Range("K61").Value = 1076.4762
$ws.Range("L61").Value = 1716.1
$ws.Range("M61").Value = -864.4762000000001
$ws.Range("N61").Value = -2140.1
$ws.Range("H116").Value = 1093.3334
$ws.Range("I116").Value = 997.7778
$ws.Range("J116").Value = 1666.6666
$ws.Range("K116").Value = 997.7778
$ws.Range("L116").Value = 1666.6666
$ws.Range("M116").Value = 1296.2222
$ws.Range("N116").Value = -6254.6666
$ws.Range("H136").Value = 1282.8064
$ws.Range("I136").Value = 1076.4762
$ws.Range("J136").Value = 1716.1
$ws.Range("K136").Value = 3229.4286
$ws.Range("L136").Value = 5148.299999999999
$ws.Range("M136").Value = -679.4286000000002
$ws.Range("N136").Value = -10248.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1093.3334
$ws.Range("I3").Value = 997.7778
$ws.Range("J3").Value = 1666.6666
$ws.Range("K3").Value = 997.7778
$ws.Range("L3").Value = 1666.6666
$ws.Range("M3").Value = -883.7778
$ws.Range("N3").Value = -1894.6666
$ws.Range("H9").Value = 48000
$ws.Range("J9").Value = 48000
$ws.Range("L9").Value = 48000
$ws.Range("N9").Value = -48336
$ws.Range("H44").Value = 48000
$ws.Range("J44").Value = 48000
$ws.Range("L44").Value = 48000
$ws.Range("N44").Value = -48994
$ws.Range("H45").Value = 59995
$ws.Range("J45").Value = 59995
$ws.Range("L45").Value = 59995
$ws.Range("N45").Value = -61611
$ws.Range("H80").Value = 699.7059
$ws.Range("J80").Value = 786.86664
$ws.Range("L80").Value = 786.86664
$ws.Range("N80").Value = -2782.86664
$ws.Range("H83").Value = 699.7059
$ws.Range("J83").Value = 786.86664
$ws.Range("L83").Value = 3934.3332
$ws.Range("N83").Value = -13918.3332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H58").Value = 1651.6666
$ws.Range("I58").Value = 1651
$ws.Range("J58").Value = 1651.8
$ws.Range("K58").Value = 1651
$ws.Range("L58").Value = 1651.8
$ws.Range("M58").Value = -1448
$ws.Range("N58").Value = -2057.8
$ws.Range("H136").Value = 1651.6666
$ws.Range("I136").Value = 1651
$ws.Range("J136").Value = 1651.8
$ws.Range("K136").Value = 4953
$ws.Range("L136").Value = 4955.4
$ws.Range("M136").Value = -2403
$ws.Range("N136").Value = -10055.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 43385.13
$ws.Range("J131").Value = 52201.82
$ws.Range("L131").Value = 156605.46
$ws.Range("N131").Value = -166685.46

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2016.7826
$ws.Range("I40").Value = 1811.625
$ws.Range("J40").Value = 2485.7144
$ws.Range("K40").Value = 1811.625
$ws.Range("L40").Value = 2485.7144
$ws.Range("M40").Value = -1675.625
$ws.Range("N40").Value = -2757.7144
$ws.Range("H68").Value = 9980.6
$ws.Range("I68").Value = 20050
$ws.Range("J68").Value = 3267.6667
$ws.Range("K68").Value = 20050
$ws.Range("L68").Value = 3267.6667
$ws.Range("M68").Value = -19301
$ws.Range("N68").Value = -4765.6667
$ws.Range("H71").Value = 9980.6
$ws.Range("I71").Value = 20050
$ws.Range("J71").Value = 3267.6667
$ws.Range("K71").Value = 100250
$ws.Range("L71").Value = 16338.3335
$ws.Range("M71").Value = -96506
$ws.Range("N71").Value = -23826.3335
$ws.Range("H136").Value = 87426
$ws.Range("I136").Value = 202762.4
$ws.Range("J136").Value = 5042.857
$ws.Range("K136").Value = 608287.2
$ws.Range("L136").Value = 15128.571
$ws.Range("M136").Value = -605737.2
$ws.Range("N136").Value = -20228.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1655.3334
$ws.Range("J8").Value = 2000
$ws.Range("L8").Value = 2000
$ws.Range("N8").Value = -2280
$ws.Range("H44").Value = 22000
$ws.Range("J44").Value = 22000
$ws.Range("L44").Value = 22000
$ws.Range("N44").Value = -23108
$ws.Range("H81").Value = 9286.875
$ws.Range("I81").Value = 28550.25
$ws.Range("J81").Value = 2865.75
$ws.Range("K81").Value = 57100.5
$ws.Range("L81").Value = 5731.5
$ws.Range("M81").Value = -56039.5
$ws.Range("N81").Value = -7853.5
$ws.Range("H84").Value = 9286.875
$ws.Range("I84").Value = 28550.25
$ws.Range("J84").Value = 2865.75
$ws.Range("K84").Value = 285502.5
$ws.Range("L84").Value = 28657.5
$ws.Range("M84").Value = -280198.5
$ws.Range("N84").Value = -39265.5
